$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ALC.Range("H33").Value = 256.52112
$ws_ALC.Range("I33").Value = 204.55072
$ws_ALC.Range("K33").Value = 204.55072
$ws_ALC.Range("M33").Value = 24.44927999999999
$ws_ALC.Range("H88").Value = 4244.5
$ws_ALC.Range("I88").Value = 5497.3
$ws_ALC.Range("J88").Value = 2678.5
$ws_ALC.Range("K88").Value = 5497.3
$ws_ALC.Range("L88").Value = 2678.5
$ws_ALC.Range("M88").Value = -5091.3
$ws_ALC.Range("N88").Value = -3490.5
$ws_ALC.Range("H91").Value = 4244.5
$ws_ALC.Range("I91").Value = 5497.3
$ws_ALC.Range("J91").Value = 2678.5
$ws_ALC.Range("K91").Value = 5497.3
$ws_ALC.Range("L91").Value = 2678.5
$ws_ALC.Range("M91").Value = -4093.3
$ws_ALC.Range("N91").Value = -5486.5
$ws_ALC.Range("H107").Value = 7429.467
$ws_ALC.Range("I107").Value = 9000.166999999999
$ws_ALC.Range("J107").Value = 1146.6666
$ws_ALC.Range("K107").Value = 9000.166999999999
$ws_ALC.Range("L107").Value = 1146.6666
$ws_ALC.Range("M107").Value = -7080.166999999999
$ws_ALC.Range("N107").Value = -4986.6666
$ws_ALC.Range("H116").Value = 2991.6667
$ws_ALC.Range("I116").Value = 2233.3333
$ws_ALC.Range("J116").Value = 3750
$ws_ALC.Range("K116").Value = 2233.3333
$ws_ALC.Range("L116").Value = 3750
$ws_ALC.Range("M116").Value = 1208.6667
$ws_ALC.Range("N116").Value = -10634
$ws_ALC.Range("H132").Value = 23584.762
$ws_ALC.Range("I132").Value = 2947.861
$ws_ALC.Range("J132").Value = 147406.17
$ws_ALC.Range("K132").Value = 8843.582999999999
$ws_ALC.Range("L132").Value = 442218.51
$ws_ALC.Range("M132").Value = -6313.582999999999
$ws_ALC.Range("N132").Value = -447278.51
$ws_ALC.Range("H135").Value = 15152223
$ws_ALC.Range("I135").Value = 738.1786
$ws_ALC.Range("J135").Value = 100000536
$ws_ALC.Range("K135").Value = 6643.6074
$ws_ALC.Range("L135").Value = 900004824
$ws_ALC.Range("M135").Value = -4108.6074
$ws_ALC.Range("N135").Value = -900009894
$ws_ALC.Range("H137").Value = 2082897.2
$ws_ALC.Range("I137").Value = 5495330
$ws_ALC.Range("K137").Value = 16485990
$ws_ALC.Range("M137").Value = -16483440
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_ARM.Range("H2").Value = 1769.2941
$ws_ARM.Range("I2").Value = 1791.3572
$ws_ARM.Range("J2").Value = 1666.3334
$ws_ARM.Range("K2").Value = 1791.3572
$ws_ARM.Range("L2").Value = 1666.3334
$ws_ARM.Range("M2").Value = -1678.3572
$ws_ARM.Range("N2").Value = -1892.3334
$ws_ARM.Range("H45").Value = 1724.125
$ws_ARM.Range("I45").Value = 1662.6923
$ws_ARM.Range("J45").Value = 1796.7273
$ws_ARM.Range("K45").Value = 1662.6923
$ws_ARM.Range("L45").Value = 1796.7273
$ws_ARM.Range("M45").Value = -1285.6923
$ws_ARM.Range("N45").Value = -2550.7273
$ws_ARM.Range("H61").Value = 1636.5358
$ws_ARM.Range("I61").Value = 1057.7391
$ws_ARM.Range("J61").Value = 4299
$ws_ARM.Range("K61").Value = 1057.7391
$ws_ARM.Range("L61").Value = 4299
$ws_ARM.Range("M61").Value = -845.7391
$ws_ARM.Range("N61").Value = -4723
$ws_ARM.Range("H74").Value = 1757.7234
$ws_ARM.Range("I74").Value = 1411.8611
$ws_ARM.Range("K74").Value = 1411.8611
$ws_ARM.Range("M74").Value = -537.8611000000001
$ws_ARM.Range("H77").Value = 1757.7234
$ws_ARM.Range("I77").Value = 1411.8611
$ws_ARM.Range("K77").Value = 7059.3055
$ws_ARM.Range("M77").Value = -2691.3055
$ws_ARM.Range("H116").Value = 1769.2941
$ws_ARM.Range("I116").Value = 1791.3572
$ws_ARM.Range("J116").Value = 1666.3334
$ws_ARM.Range("K116").Value = 1791.3572
$ws_ARM.Range("L116").Value = 1666.3334
$ws_ARM.Range("M116").Value = 502.6428000000001
$ws_ARM.Range("N116").Value = -6254.3334
$ws_ARM.Range("H132").Value = 2046.8
$ws_ARM.Range("I132").Value = 1129.3334
$ws_ARM.Range("K132").Value = 3388.0002
$ws_ARM.Range("M132").Value = -858.0001999999999
$ws_ARM.Range("H136").Value = 1636.5358
$ws_ARM.Range("I136").Value = 1057.7391
$ws_ARM.Range("J136").Value = 4299
$ws_ARM.Range("K136").Value = 3173.2173
$ws_ARM.Range("L136").Value = 12897
$ws_ARM.Range("M136").Value = -623.2173000000003
$ws_ARM.Range("N136").Value = -17997
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_BSM.Range("H3").Value = 1769.2941
$ws_BSM.Range("I3").Value = 1791.3572
$ws_BSM.Range("J3").Value = 1666.3334
$ws_BSM.Range("K3").Value = 1791.3572
$ws_BSM.Range("L3").Value = 1666.3334
$ws_BSM.Range("M3").Value = -1677.3572
$ws_BSM.Range("N3").Value = -1894.3334
$ws_BSM.Range("H134").Value = 2198.277
$ws_BSM.Range("I134").Value = 1380.5098
$ws_BSM.Range("J134").Value = 3501.5938
$ws_BSM.Range("K134").Value = 4141.5294
$ws_BSM.Range("L134").Value = 10504.7814
$ws_BSM.Range("M134").Value = -1606.5294
$ws_BSM.Range("N134").Value = -15574.7814
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CRP.Range("H31").Value = 4237.4
$ws_CRP.Range("I31").Value = 1325.1852
$ws_CRP.Range("K31").Value = 1325.1852
$ws_CRP.Range("M31").Value = -1030.1852
$ws_CRP.Range("H34").Value = 4237.4
$ws_CRP.Range("I34").Value = 1325.1852
$ws_CRP.Range("K34").Value = 1325.1852
$ws_CRP.Range("M34").Value = -1123.1852
$ws_CRP.Range("H58").Value = 1722.5122
$ws_CRP.Range("I58").Value = 1140
$ws_CRP.Range("K58").Value = 1140
$ws_CRP.Range("M58").Value = -937
$ws_CRP.Range("H132").Value = 75665.94500000001
$ws_CRP.Range("I132").Value = 1319.4
$ws_CRP.Range("J132").Value = 158273.22
$ws_CRP.Range("K132").Value = 3958.2
$ws_CRP.Range("L132").Value = 474819.66
$ws_CRP.Range("M132").Value = -1428.2
$ws_CRP.Range("N132").Value = -479879.66
$ws_CRP.Range("H134").Value = 1192716.1
$ws_CRP.Range("I134").Value = 1566123.6
$ws_CRP.Range("J134").Value = 352549.25
$ws_CRP.Range("K134").Value = 4698370.800000001
$ws_CRP.Range("L134").Value = 1057647.75
$ws_CRP.Range("M134").Value = -4695835.800000001
$ws_CRP.Range("N134").Value = -1062717.75
$ws_CRP.Range("H136").Value = 1722.5122
$ws_CRP.Range("I136").Value = 1140
$ws_CRP.Range("K136").Value = 3420
$ws_CRP.Range("M136").Value = -870
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_CUL.Range("H38").Value = 105438.69
$ws_CUL.Range("I38").Value = 104.09091
$ws_CUL.Range("J38").Value = 250273.75
$ws_CUL.Range("K38").Value = 312.27273
$ws_CUL.Range("L38").Value = 750821.25
$ws_CUL.Range("M38").Value = 34.72727000000003
$ws_CUL.Range("N38").Value = -751515.25
$ws_CUL.Range("H46").Value = 1560
$ws_CUL.Range("I46").Value = 0
$ws_CUL.Range("J46").Value = 1560
$ws_CUL.Range("K46").Value = 0
$ws_CUL.Range("L46").Value = 4680
$ws_CUL.Range("M46").ClearContents()
$ws_CUL.Range("N46").Value = -4862
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_GSM.Range("H132").Value = 40005064
$ws_GSM.Range("I132").Value = 76928650
$ws_GSM.Range("J132").Value = 4517.1665
$ws_GSM.Range("K132").Value = 230785950
$ws_GSM.Range("L132").Value = 13551.4995
$ws_GSM.Range("M132").Value = -230783420
$ws_GSM.Range("N132").Value = -18611.4995
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_LTW.Range("H7").Value = 2939.75
$ws_LTW.Range("I7").Value = 2278.75
$ws_LTW.Range("J7").Value = 3931.25
$ws_LTW.Range("K7").Value = 2278.75
$ws_LTW.Range("L7").Value = 3931.25
$ws_LTW.Range("M7").Value = -2166.75
$ws_LTW.Range("N7").Value = -4155.25
$ws_LTW.Range("H61").Value = 2874.5881
$ws_LTW.Range("I61").Value = 3135.75
$ws_LTW.Range("J61").Value = 2642.4443
$ws_LTW.Range("K61").Value = 3135.75
$ws_LTW.Range("L61").Value = 2642.4443
$ws_LTW.Range("M61").Value = -2933.75
$ws_LTW.Range("N61").Value = -3046.4443
$ws_LTW.Range("H94").Value = 83000
$ws_LTW.Range("J94").Value = 83000
$ws_LTW.Range("L94").Value = 83000
$ws_LTW.Range("N94").Value = -84352
$ws_LTW.Range("H113").Value = 2874.5881
$ws_LTW.Range("I113").Value = 3135.75
$ws_LTW.Range("J113").Value = 2642.4443
$ws_LTW.Range("K113").Value = 3135.75
$ws_LTW.Range("L113").Value = 2642.4443
$ws_LTW.Range("M113").Value = -965.75
$ws_LTW.Range("N113").Value = -6982.4443
$ws_LTW.Range("H126").Value = 2939.75
$ws_LTW.Range("I126").Value = 2278.75
$ws_LTW.Range("J126").Value = 3931.25
$ws_LTW.Range("K126").Value = 6836.25
$ws_LTW.Range("L126").Value = 11793.75
$ws_LTW.Range("M126").Value = -4366.25
$ws_LTW.Range("N126").Value = -16733.75
$ws_LTW.Range("H132").Value = 3033.7188
$ws_LTW.Range("I132").Value = 2179.6
$ws_LTW.Range("J132").Value = 4457.25
$ws_LTW.Range("K132").Value = 6538.799999999999
$ws_LTW.Range("L132").Value = 13371.75
$ws_LTW.Range("M132").Value = -4008.799999999999
$ws_LTW.Range("N132").Value = -18431.75
$ws_LTW.Range("H136").Value = 1980.2413
$ws_LTW.Range("I136").Value = 1541.4445
$ws_LTW.Range("J136").Value = 2698.2727
$ws_LTW.Range("K136").Value = 4624.333500000001
$ws_LTW.Range("L136").Value = 8094.8181
$ws_LTW.Range("M136").Value = -2074.333500000001
$ws_LTW.Range("N136").Value = -13194.8181
$ws_WVR = $wb.Worksheets.Item("WVR")
$ws_WVR.Range("H104").Value = 43010
$ws_WVR.Range("J104").Value = 43010
$ws_WVR.Range("L104").Value = 43010
$ws_WVR.Range("N104").Value = -49998
$ws_WVR.Range("H132").Value = 1176647.4
$ws_WVR.Range("J132").Value = 3528
$ws_WVR.Range("L132").Value = 10584
$ws_WVR.Range("N132").Value = -15644
$ws_WVR.Range("H136").Value = 973629.9399999999
$ws_WVR.Range("I136").Value = 1667666.5
$ws_WVR.Range("J136").Value = 1978.7
$ws_WVR.Range("K136").Value = 5002999.5
$ws_WVR.Range("L136").Value = 5936.1
$ws_WVR.Range("M136").Value = -5000449.5
$ws_WVR.Range("N136").Value = -11036.1
